# Rebuild the "Estado de Cuenta" detail table (rows 16-34) on Hoja1.
# Previous account-statement periods are removed and replaced with the
# new set of periods per worker (see commit message: "Elimina EC
# anteriores y se agregan nuevos, se modifica base de datos").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# IGNACIO TORRES BALSEIRO (CC 73203467) - rows 16-27
$ws.Range("B16:B27").Value = "CC"
$ws.Range("C16:C27").Value = "73203467"
$ws.Range("D16:D27").Value = "IGNACIO TORRES BALSEIRO"

$ws.Range("E16").Value = "2211"
$ws.Range("F16").Value = 25439
$ws.Range("G16").Value = 908526

$ws.Range("E17").Value = "2210"
$ws.Range("F17").Value = 36341
$ws.Range("G17").Value = 908526

$ws.Range("E18").Value = "2209"
$ws.Range("F18").Value = 36341
$ws.Range("G18").Value = 908526

$ws.Range("E19").Value = "2208"
$ws.Range("F19").Value = 36341
$ws.Range("G19").Value = 908526

$ws.Range("E20").Value = "2207"
$ws.Range("F20").Value = 36341
$ws.Range("G20").Value = 908526

$ws.Range("E21").Value = "2206"
$ws.Range("F21").Value = 36341
$ws.Range("G21").Value = 908526

$ws.Range("E22").Value = "2205"
$ws.Range("F22").Value = 36341
$ws.Range("G22").Value = 908526

$ws.Range("E23").Value = "2204"
$ws.Range("F23").Value = 36341
$ws.Range("G23").Value = 908526

$ws.Range("E24").Value = "2203"
$ws.Range("F24").Value = 36341
$ws.Range("G24").Value = 908526

$ws.Range("E25").Value = "2202"
$ws.Range("F25").Value = 36341
$ws.Range("G25").Value = 908526

$ws.Range("E26").Value = "2201"
$ws.Range("F26").Value = 36341
$ws.Range("G26").Value = 908526

$ws.Range("E27").Value = "2112"
$ws.Range("F27").Value = 10902
$ws.Range("G27").Value = 908526

# SANDY JAVIER DURANGO PEREGRINO (CC 1143326442) - rows 28-34
$ws.Range("B28:B34").Value = "CC"
$ws.Range("C28:C34").Value = "1143326442"
$ws.Range("D28:D34").Value = "SANDY JAVIER DURANGO PEREGRINO"

$ws.Range("E28").Value = "2211"
$ws.Range("F28").Value = 33419
$ws.Range("G28").Value = 1193546

$ws.Range("E29").Value = "2210"
$ws.Range("F29").Value = 47742
$ws.Range("G29").Value = 1193546

$ws.Range("E30").Value = "2209"
$ws.Range("F30").Value = 47742
$ws.Range("G30").Value = 1193546

$ws.Range("E31").Value = "2208"
$ws.Range("F31").Value = 47742
$ws.Range("G31").Value = 1193546

$ws.Range("E32").Value = "2207"
$ws.Range("F32").Value = 47742
$ws.Range("G32").Value = 1193546

$ws.Range("E33").Value = "2206"
$ws.Range("F33").Value = 47742
$ws.Range("G33").Value = 1193546

$ws.Range("E34").Value = "2205"
$ws.Range("F34").Value = 47742
$ws.Range("G34").Value = 1193546
